# Rename the observed-variable headers from the Soybean model names to the
# Mungbean model names (the sheet was repurposed for Mungbean, not Soybean).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Observed")

$ws.Range("C1").Value = "Mungbean.Phenology.CurrentStageName"
$ws.Range("D1").Value = "Mungbean.Phenology.MaturityDAS"
$ws.Range("E1").Value = "Mungbean.AboveGround.Wt"
$ws.Range("F1").Value = "Mungbean.AboveGround.Wterror"
$ws.Range("G1").Value = "Mungbean.Grain.Wt"
$ws.Range("H1").Value = "Mungbean.Grain.Wterror"
$ws.Range("K1").Value = "Mungbean.Grain.HarvestIndex"
$ws.Range("L1").Value = "Mungbean.Grain.HarvestIndexerror"

# Move the active selection, matching the author's last cursor position.
$ws.Range("K18").Select()
